$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-25 Wednesday" "2026-02-26 Thursday"

Replace-Text "65÷7=" "84÷7="
Replace-Text "70÷2=" "65÷5="
Replace-Text "85÷5=" "54÷9="
Replace-Text "49÷5=" "53÷7="
Replace-Text "44÷5=" "53÷4="
Replace-Text "87÷8=" "42÷7="
Replace-Text "88÷3=" "80÷9="
Replace-Text "14÷5=" "71÷7="
Replace-Text "36÷2=" "65÷3="
Replace-Text "70÷9=" "40÷7="
Replace-Text "50÷6=" "33÷6="
Replace-Text "98÷2=" "65÷5="
Replace-Text "81÷3=" "89÷6="
Replace-Text "14÷9=" "82÷6="
Replace-Text "74÷6=" "54÷2="
Replace-Text "92÷3=" "68÷5="
Replace-Text "85÷9=" "11÷6="
Replace-Text "53÷9=" "76÷7="
Replace-Text "67÷9=" "90÷2="
Replace-Text "81÷6=" "74÷2="
Replace-Text "11÷3=" "90÷4="
Replace-Text "21÷2=" "91÷7="
Replace-Text "43÷6=" "28÷6="
Replace-Text "16÷5=" "60÷4="
Replace-Text "34÷6=" "48÷7="
